$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value2 = "Datos actualizados a 22 de Marzo de 2020 a las 01:16"

# Rows 53 and 54 fully swap (Huelva <-> Huesca), including their "Casos activos" values
$a53 = $ws.Range("A53").Value2
$b53 = $ws.Range("B53").Value2
$c53 = $ws.Range("C53").Value2
$d53 = $ws.Range("D53").Value2
$e53 = $ws.Range("E53").Value2

$a54 = $ws.Range("A54").Value2
$b54 = $ws.Range("B54").Value2
$c54 = $ws.Range("C54").Value2
$d54 = $ws.Range("D54").Value2
$e54 = $ws.Range("E54").Value2

$ws.Range("A53").Value2 = $a54
$ws.Range("B53").Value2 = $b54
$ws.Range("C53").Value2 = $c54
$ws.Range("D53").Value2 = $d54
$ws.Range("E53").Value2 = $e54

$ws.Range("A54").Value2 = $a53
$ws.Range("B54").Value2 = $b53
$ws.Range("C54").Value2 = $c53
$ws.Range("D54").Value2 = $d53
$ws.Range("E54").Value2 = $e53
